# Updated symbol list on Sun Jan 15 04:59:01 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# several coins on the active sheet. Values are written with a leading
# apostrophe so Excel keeps them as literal text (matching the workbook's
# existing convention of storing these figures as text strings) instead of
# auto-converting them to numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'298.44"
$ws.Range("E2").Value = "'-2.70%"

$ws.Range("E3").Value = "'-1.58%"

$ws.Range("D4").Value = "'5.059"
$ws.Range("E4").Value = "'-5.25%"

$ws.Range("D5").Value = "'0.07515"
$ws.Range("E5").Value = "'0.81%"

$ws.Range("D6").Value = "'7.793"
$ws.Range("E6").Value = "'0.64%"

$ws.Range("D7").Value = "'1.715"
$ws.Range("E7").Value = "'8.02%"

$ws.Range("E8").Value = "'2.70%"

$ws.Range("D9").Value = "'0.9263"
$ws.Range("E9").Value = "'0.64%"

$ws.Range("D10").Value = "'0.1706"
$ws.Range("E10").Value = "'1.85%"

$ws.Range("D11").Value = "'0.07429"
$ws.Range("E11").Value = "'-3.07%"

$ws.Range("D12").Value = "'0.07976"
$ws.Range("E12").Value = "'-0.16%"

$ws.Range("D13").Value = "'0.03046"
$ws.Range("E13").Value = "'-1.04%"

$ws.Range("D14").Value = "'0.09889"
$ws.Range("E14").Value = "'0.26%"

$ws.Range("D15").Value = "'0.001499"
$ws.Range("E15").Value = "'-2.08%"

$ws.Range("D16").Value = "'0.04662"
$ws.Range("E16").Value = "'2.16%"

$ws.Range("D17").Value = "'0.006343"
$ws.Range("E17").Value = "'-1.80%"

$ws.Range("E18").Value = "'-0.10%"

$ws.Range("D19").Value = "'2.218"
$ws.Range("E19").Value = "'-1.15%"

$ws.Range("E20").Value = "'0.70%"

$ws.Range("E21").Value = "'1.39%"

$ws.Range("D22").Value = "'4.555"
$ws.Range("E22").Value = "'7.95%"

$ws.Range("D24").Value = "'0.001215"
$ws.Range("E24").Value = "'0.09%"

$ws.Range("D25").Value = "'0.004428"
$ws.Range("E25").Value = "'-2.25%"

$ws.Range("E26").Value = "'19.81%"

$ws.Range("E27").Value = "'5.99%"

$ws.Range("D39").Value = "'0.01668"
$ws.Range("E39").Value = "'1.49%"

$ws.Range("D40").Value = "'0.04546"
$ws.Range("E40").Value = "'0.54%"

$ws.Range("D41").Value = "'0.007016"
$ws.Range("E41").Value = "'-5.52%"

$ws.Range("D42").Value = "'0.1327"
$ws.Range("E42").Value = "'-2.88%"

$ws.Range("D43").Value = "'0.002060"

$ws.Range("D44").Value = "'0.01281"
$ws.Range("E44").Value = "'-6.43%"

$ws.Range("D45").Value = "'0.00006089"
$ws.Range("E45").Value = "'-0.62%"

$ws.Range("D46").Value = "'0.7116"
$ws.Range("E46").Value = "'-62.40%"

$ws.Range("E47").Value = "'-5.64%"
